$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 869.23254
$ws.Range("J17").Value = 887.8095
$ws.Range("L17").Value = 2663.4285
$ws.Range("N17").Value = -2999.4285

$ws.Range("H41").Value = 20835008
$ws.Range("I41").Value = 30304320
$ws.Range("J41").Value = 2518
$ws.Range("K41").Value = 30304320
$ws.Range("L41").Value = 2518
$ws.Range("M41").Value = -30303880
$ws.Range("N41").Value = -3398.3333

$ws.Range("H137").Value = 1487.0869
$ws.Range("I137").Value = 900.06665
$ws.Range("J137").Value = 2587.75
$ws.Range("K137").Value = 2700.19995
$ws.Range("L137").Value = 7763.25
$ws.Range("M137").Value = -150.1999500000002
$ws.Range("N137").Value = -12863.25

$ws.Range("H138").Value = 1886.8265
$ws.Range("I138").Value = 729.6111
$ws.Range("J138").Value = 2147.2
$ws.Range("K138").Value = 2188.8333
$ws.Range("L138").Value = 6441.599999999999
$ws.Range("M138").Value = 2951.1667
$ws.Range("N138").Value = -16721.6

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1444.2307
$ws.Range("I61").Value = 1052.375
$ws.Range("J61").Value = 2071.2
$ws.Range("K61").Value = 1052.375
$ws.Range("L61").Value = 2071.2
$ws.Range("M61").Value = -840.375
$ws.Range("N61").Value = -2495.2

$ws.Range("H74").Value = 1292.0588
$ws.Range("I74").Value = 1119.75
$ws.Range("J74").Value = 1705.6
$ws.Range("K74").Value = 1119.75
$ws.Range("L74").Value = 1705.6
$ws.Range("M74").Value = -245.75
$ws.Range("N74").Value = -3453.6

$ws.Range("H77").Value = 1292.0588
$ws.Range("I77").Value = 1119.75
$ws.Range("J77").Value = 1705.6
$ws.Range("K77").Value = 5598.75
$ws.Range("L77").Value = 8528
$ws.Range("M77").Value = -1230.75
$ws.Range("N77").Value = -17264

$ws.Range("H94").Value = 50000
$ws.Range("J94").Value = 50000
$ws.Range("L94").Value = 50000
$ws.Range("N94").Value = -51802

$ws.Range("H123").Value = 68071.5
$ws.Range("J123").Value = 68071.5
$ws.Range("L123").Value = 68071.5
$ws.Range("N123").Value = -77871.5

$ws.Range("H136").Value = 1444.2307
$ws.Range("I136").Value = 1052.375
$ws.Range("J136").Value = 2071.2
$ws.Range("K136").Value = 3157.125
$ws.Range("L136").Value = 6213.599999999999
$ws.Range("M136").Value = -607.125
$ws.Range("N136").Value = -11313.6

$ws.Range("H138").Value = 30216
$ws.Range("J138").Value = 30216
$ws.Range("L138").Value = 30216
$ws.Range("N138").Value = -40496

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H26").Value = 19490.334
$ws.Range("I26").Value = 19490.334
$ws.Range("K26").Value = 19490.334
$ws.Range("M26").Value = -19198.334

$ws.Range("H96").Value = 9571.200000000001
$ws.Range("I96").Value = 4214
$ws.Range("J96").Value = 31000
$ws.Range("K96").Value = 4214
$ws.Range("L96").Value = 31000
$ws.Range("M96").Value = -1468
$ws.Range("N96").Value = -36492

$ws.Range("H134").Value = 4176.1943
$ws.Range("I134").Value = 1418.1333
$ws.Range("J134").Value = 17966.5
$ws.Range("K134").Value = 4254.3999
$ws.Range("L134").Value = 53899.5
$ws.Range("M134").Value = -1719.3999
$ws.Range("N134").Value = -58969.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 3000
$ws.Range("J4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("N4").ClearContents()

$ws.Range("H31").Value = 1560.4814
$ws.Range("I31").Value = 1006
$ws.Range("J31").Value = 2075.3572
$ws.Range("K31").Value = 1006
$ws.Range("L31").Value = 2075.3572
$ws.Range("M31").Value = -711
$ws.Range("N31").Value = -2665.3572

$ws.Range("H34").Value = 1560.4814
$ws.Range("I34").Value = 1006
$ws.Range("J34").Value = 2075.3572
$ws.Range("K34").Value = 1006
$ws.Range("L34").Value = 2075.3572
$ws.Range("M34").Value = -804
$ws.Range("N34").Value = -2479.3572

$ws.Range("H95").Value = 24204.8
$ws.Range("J95").Value = 24204.8
$ws.Range("L95").Value = 24204.8
$ws.Range("N95").Value = -29696.8

$ws.Range("H134").Value = 9525141
$ws.Range("I134").Value = 12346957
$ws.Range("J134").Value = 1510.875
$ws.Range("K134").Value = 37040871
$ws.Range("L134").Value = 4532.625
$ws.Range("M134").Value = -37038336
$ws.Range("N134").Value = -9602.625

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 450322.34
$ws.Range("J4").Value = 577530.1
$ws.Range("L4").Value = 1732590.3
$ws.Range("N4").Value = -1732814.3

$ws.Range("H68").Value = 1728.3846
$ws.Range("I68").Value = 752.6
$ws.Range("J68").Value = 3059
$ws.Range("K68").Value = 2257.8
$ws.Range("L68").Value = 9177
$ws.Range("M68").Value = -1446.8
$ws.Range("N68").Value = -10799

$ws.Range("H71").Value = 1728.3846
$ws.Range("I71").Value = 752.6
$ws.Range("J71").Value = 3059
$ws.Range("K71").Value = 6773.400000000001
$ws.Range("L71").Value = 27531
$ws.Range("M71").Value = -2717.400000000001
$ws.Range("N71").Value = -35643

$ws.Range("H75").Value = 3000
$ws.Range("J75").Value = 3000
$ws.Range("L75").Value = 9000
$ws.Range("N75").Value = -10996

$ws.Range("H78").Value = 3000
$ws.Range("J78").Value = 3000
$ws.Range("L78").Value = 27000
$ws.Range("N78").Value = -36984

$ws.Range("H81").Value = 2483.5454
$ws.Range("I81").Value = 1298.8334
$ws.Range("J81").Value = 2927.8125
$ws.Range("K81").Value = 3896.5002
$ws.Range("L81").Value = 8783.4375
$ws.Range("M81").Value = -2773.5002
$ws.Range("N81").Value = -11029.4375

$ws.Range("H84").Value = 2483.5454
$ws.Range("I84").Value = 1298.8334
$ws.Range("J84").Value = 2927.8125
$ws.Range("K84").Value = 11689.5006
$ws.Range("L84").Value = 26350.3125
$ws.Range("M84").Value = -6073.500599999999
$ws.Range("N84").Value = -37582.3125

$ws.Range("H99").Value = 1865.5834
$ws.Range("J99").Value = 2564
$ws.Range("L99").Value = 7692
$ws.Range("N99").Value = -12184

$ws.Range("H102").Value = 3889
$ws.Range("J102").Value = 3889
$ws.Range("L102").Value = 11667
$ws.Range("N102").Value = -16535

$ws.Range("H123").Value = 1074.75
$ws.Range("J123").Value = 3999
$ws.Range("L123").Value = 11997
$ws.Range("N123").Value = -16897

$ws.Range("H131").Value = 23257310
$ws.Range("J131").Value = 1602.6757
$ws.Range("L131").Value = 4808.0271
$ws.Range("N131").Value = -14888.0271

$ws.Range("H139").Value = 2616.5715
$ws.Range("J139").Value = 1999.1538
$ws.Range("L139").Value = 5997.4614
$ws.Range("N139").Value = -16277.4614

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H7").Value = 17859.715
$ws.Range("J7").Value = 17859.715
$ws.Range("L7").Value = 17859.715
$ws.Range("N7").Value = -18083.715

$ws.Range("H8").Value = 17859.715
$ws.Range("J8").Value = 17859.715
$ws.Range("L8").Value = 17859.715
$ws.Range("N8").Value = -18137.715

$ws.Range("H11").Value = 6037500
$ws.Range("I11").Value = 6037500
$ws.Range("K11").Value = 6037500
$ws.Range("M11").Value = -6037361

$ws.Range("H102").Value = 2001.15
$ws.Range("I102").Value = 1973.0714
$ws.Range("J102").Value = 2066.6667
$ws.Range("K102").Value = 1973.0714
$ws.Range("L102").Value = 2066.6667
$ws.Range("M102").Value = -351.0714
$ws.Range("N102").Value = -5310.6667

$ws.Range("H133").Value = 38649.25
$ws.Range("J133").Value = 38649.25
$ws.Range("L133").Value = 38649.25
$ws.Range("N133").Value = -48769.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 3008.8
$ws.Range("J2").Value = 3123.111
$ws.Range("L2").Value = 3123.111
$ws.Range("N2").Value = -3347.111

$ws.Range("H16").Value = 1089.2858
$ws.Range("I16").Value = 863.2353000000001
$ws.Range("K16").Value = 863.2353000000001
$ws.Range("M16").Value = -693.2353000000001

$ws.Range("H74").Value = 15000
$ws.Range("I74").Value = 15000
$ws.Range("K74").Value = 15000
$ws.Range("M74").Value = -14002

$ws.Range("H77").Value = 15000
$ws.Range("I77").Value = 15000
$ws.Range("K77").Value = 45000
$ws.Range("M77").Value = -40008

$ws.Range("H122").Value = 9619110
$ws.Range("I122").Value = 16670128
$ws.Range("K122").Value = 50010384
$ws.Range("M122").Value = -50007934

$ws.Range("H132").Value = 40510.348
$ws.Range("I132").Value = 1684.0667
$ws.Range("J132").Value = 93455.27
$ws.Range("K132").Value = 5052.2001
$ws.Range("L132").Value = 280365.81
$ws.Range("M132").Value = -2522.2001
$ws.Range("N132").Value = -285425.81

$ws.Range("H136").Value = 8066.4
$ws.Range("I136").Value = 12153.667
$ws.Range("K136").Value = 36461.001
$ws.Range("M136").Value = -33911.001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H104").Value = 22000
$ws.Range("J104").Value = 22000
$ws.Range("L104").Value = 22000
$ws.Range("N104").Value = -28988
